# Add two new columns (I: I0, J: IF) to the sheet, mirroring column H's
# header style, and populate the per-row numeric data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells (I1, J1) ---------------------------------------------
# Copy H1's formatting (bold font, border, centered/top alignment) onto
# the two new header cells so they reuse the same style as the rest of
# row 1, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data cells (I2:J77) ------------------------------------------------
$iValues = @(9,10,9,9,9,9,9,9,9,8,8,7,8,7,8,8,8,8,10,8,6,6,7,8,7,6,6,6,8,8,7,10,11,7,8,8,9,8,9,9,9,8,8,9,9,11,8,8,7,8,9,8,7,9,6,9,8,9,6,8,9,8,9,8,6,8,8,6,6,7,6,7,7,6,6,6)
$jValues = @(9,10,9,9,9,9,9,9,9,8,8,8,8,8,8,8,8,8,10,8,7,6,7,8,7,6,6,6,8,8,8,10,12,7,8,8,9,8,9,9,9,8,8,9,9,11,8,8,7,8,9,8,8,9,7,9,8,9,6,8,9,8,9,8,6,8,8,6,6,7,6,7,7,6,6,6)

for ($r = 0; $r -lt $iValues.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
